# Delete the empty placeholder rows 15:18 on the first sheet ("بورد اصلی").
# This shifts every row at/after 19 up by 4, collapsing the gap left by the
# two blank rows (16 & 17) that used to sit between the "توضيحات" block and
# the "عضويت" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows("15:18").Delete() | Out-Null

# The plain (non-table) AutoFilter range and the two dependent defined
# names do not auto-shrink when rows are deleted, so re-anchor them to the
# new used range (was A1:E47 / D1:E47, now 4 rows shorter).
if ($ws.AutoFilterMode) { $ws.AutoFilterMode = $false }
$ws.Range("A1:E43").AutoFilter() | Out-Null

$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "='بورد اصلی'!`$A`$1:`$E`$43"
$wb.Names.Item("_xlcn.WorksheetConnection_بورداصلیD1E311").RefersTo = "='بورد اصلی'!`$D`$1:`$E`$43"

# Match the post-edit selection recorded for the sheet.
$ws.Range("A21").Select() | Out-Null
